# Applies the "Attributes" sheet addition + related workbook changes
# described by the commit "Ergänzung um Funktionalität zur Behandlung von
# Klassenattributen" (addition of functionality for handling class
# attributes).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("ClassBounds")
$ws2 = $wb.Worksheets.Item("Tabelle2")

# --- Rename Tabelle2 -> Attributes ------------------------------------
$ws2.Name = "Attributes"

# --- Populate the Attributes sheet ------------------------------------
# Header row (row 1) - same "heading" style as used on ClassBounds row 1
$ws2.Range("A1").Value = "Class"
$ws2.Range("B1").Value = "Attribute"
$ws2.Range("C1").Value = "Concrete (mandatory)"
$ws2.Range("D1").Value = "Undefined fix"
$ws2.Range("E1").Value = "Set elements fix"
$ws2.Range("F1").Value = "Domain"
$ws2.Range("G1").Value = "Minimum Number"
$ws2.Range("H1").Value = "Maximum Number"
$ws2.Range("J1").Value = "Domain"
$ws2.Range("K1").Value = "Values"
$ws2.Range("A1:K1").Style = "Überschrift 2"
$ws2.Rows.Item(1).RowHeight = 18

# Row 2 stays empty (separator row, mirrors ClassBounds row 2 height)
$ws2.Rows.Item(2).RowHeight = 15.75

# Row 3 - "name" attribute
$ws2.Range("A3").Value = "Person"
$ws2.Range("B3").Value = "name"
$ws2.Range("C3").Value = "ada->'Ada',bob->'Bob'"
$ws2.Range("F3").Value = "names"
$ws2.Range("J3").Value = "names"
$ws2.Range("K3").Value = "'Ada', 'Bob', 'Cyd', 'Dan'"

# Row 4 - "age" attribute
$ws2.Range("A4").Value = "Person"
$ws2.Range("B4").Value = "age"
$ws2.Range("C4").Value = "ada->3"
$ws2.Range("F4").Value = "ages"
$ws2.Range("J4").Value = "numbers"
$ws2.Range("K4").Value = "1,2,3,4,5,12"

# Row 5 - "luckyNumbers" attribute
$ws2.Range("A5").Value = "Person"
$ws2.Range("B5").Value = "luckyNumbers"
$ws2.Range("C5").Value = "ada->Set{1,2,3}"
$ws2.Range("F5").Value = "numbers"
$ws2.Range("J5").Value = "ages"
$ws2.Range("K5").Value = "1..100"

# --- Column widths (best-effort match of the authored widths) --------
$ws2.Columns.Item(1).ColumnWidth = 10.5703125
$ws2.Columns.Item(2).ColumnWidth = 13.7109375
$ws2.Columns.Item(3).ColumnWidth = 23.7109375
$ws2.Columns.Item(4).ColumnWidth = 14.85546875
$ws2.Columns.Item(5).ColumnWidth = 17.85546875
$ws2.Columns.Item(6).ColumnWidth = 20.140625
$ws2.Columns.Item(7).ColumnWidth = 20.5703125
$ws2.Columns.Item(8).ColumnWidth = 20.5703125
$ws2.Columns.Item(9).ColumnWidth = 10.85546875
$ws2.Columns.Item(11).ColumnWidth = 20.7109375

# --- Selection state on ClassBounds (no longer the active tab) -------
$ws1.Cells.Select()

# --- Activate Attributes sheet with its own selection -----------------
$ws2.Range("D3").Select()

Write-Output "Applied Attributes sheet edits"
